# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update several countries' case/death/recovered counters
# - Re-sort a few rows alphabetically by country name (Cabo Verde,
#   Liechtenstein, Barbados, Zambia / Bahamas, Guyana, Islas Caimanes /
#   Burundi, Islas Turcas y Caicos, Seychelles)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 23 de Abril de 2020 a las 13:52"

# Row 8
$ws.Cells.Item(8,2).Value = 150773
$ws.Cells.Item(8,3).Value = 125
$ws.Cells.Item(8,5).Value = 42155
$ws.Cells.Item(8,7).Value = 3
$ws.Cells.Item(8,8).Value = 5318

# Row 51
$ws.Cells.Item(51,5).Value = 2112
$ws.Cells.Item(51,6).Value = 60
$ws.Cells.Item(51,7).Value = 23
$ws.Cells.Item(51,8).Value = 172

# Row 69
$ws.Cells.Item(69,4).Value = 486
$ws.Cells.Item(69,5).Value = 1223

# Row 110
$ws.Cells.Item(110,4).Value = 111
$ws.Cells.Item(110,5).Value = 304

# Row 150
$ws.Cells.Item(150,1).Value = "Cabo Verde"
$ws.Cells.Item(150,2).Value = 82
$ws.Cells.Item(150,3).Value = 9
$ws.Cells.Item(150,4).Value = 1
$ws.Cells.Item(150,5).Value = 80

# Row 151
$ws.Cells.Item(151,1).Value = "Liechtenstein"
$ws.Cells.Item(151,2).Value = 81
$ws.Cells.Item(151,4).Value = 55
$ws.Cells.Item(151,5).Value = 25
$ws.Cells.Item(151,6).Value = 0
$ws.Cells.Item(151,8).Value = 1

# Row 152
$ws.Cells.Item(152,1).Value = "Barbados"
$ws.Cells.Item(152,2).Value = 76
$ws.Cells.Item(152,4).Value = 27
$ws.Cells.Item(152,5).Value = 43
$ws.Cells.Item(152,6).Value = 4
$ws.Cells.Item(152,8).Value = 6

# Row 153
$ws.Cells.Item(153,1).Value = "Zambia"
$ws.Cells.Item(153,2).Value = 74
$ws.Cells.Item(153,4).Value = 35
$ws.Cells.Item(153,5).Value = 36
$ws.Cells.Item(153,6).Value = 1
$ws.Cells.Item(153,8).Value = 3

# Row 155
$ws.Cells.Item(155,1).Value = "Bahamas"
$ws.Cells.Item(155,2).Value = 70
$ws.Cells.Item(155,3).Value = 5
$ws.Cells.Item(155,4).Value = 12
$ws.Cells.Item(155,5).Value = 49
$ws.Cells.Item(155,6).Value = 1
$ws.Cells.Item(155,8).Value = 9

# Row 156
$ws.Cells.Item(156,1).Value = "Guyana"
$ws.Cells.Item(156,2).Value = 67
$ws.Cells.Item(156,4).Value = 9
$ws.Cells.Item(156,5).Value = 51
$ws.Cells.Item(156,6).Value = 5
$ws.Cells.Item(156,8).Value = 7

# Row 157
$ws.Cells.Item(157,1).Value = "Islas Caimanes"
$ws.Cells.Item(157,2).Value = 66
$ws.Cells.Item(157,4).Value = 7
$ws.Cells.Item(157,5).Value = 58
$ws.Cells.Item(157,6).Value = 3
$ws.Cells.Item(157,8).Value = 1

# Row 197
$ws.Cells.Item(197,1).Value = "Burundi"
$ws.Cells.Item(197,4).Value = 4
$ws.Cells.Item(197,8).Value = 1

# Row 198
$ws.Cells.Item(198,1).Value = "Islas Turcas y Caicos"

# Row 199
$ws.Cells.Item(199,1).Value = "Seychelles"
$ws.Cells.Item(199,4).Value = 6
$ws.Cells.Item(199,5).Value = 5
$ws.Cells.Item(199,8).Value = 0
